$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-0.163, -0.108, 0.212, 94, 6, -78),
    @(0.397, -0.198, -0.455, -7, -30, -6),
    @(-0.308, -0.233, 0.6840000000000001, 94, 85, 42),
    @(-0.036, -0.9379999999999999, 0.119, -20, 40, -54),
    @(-0.22, -0.07199999999999999, 0.192, -146, 50, -133),
    @(0.128, 0.7, 0.145, -21, -10, 176),
    @(0.092, -0.152, -0.303, -25, -26, 38),
    @(-0.261, -0.316, 0.701, 122, -49, 97),
    @(-0.341, -0.215, 0.717, 93, -60, 14),
    @(-0.136, 0.235, 0.1, -137, -14, 19)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 1
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $row[$j]
    }
}
